$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row: "<field>_old" -> "<field>_FV2310" and
#    "<field>_new" -> "<field>_FV2404" (the "diff" column, col K, is unchanged).
$fields = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $fields.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = ($fields[$i] + "_FV2310")
}
for ($i = 0; $i -lt $fields.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = ($fields[$i] + "_FV2404")
}

# 2) Turn the used range into an Excel table ("Table1") with a header row.
$tableRange = $ws.Range("A1:U78")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

# 3) Freeze the header row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
